$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a text value to a cell while forcing text storage
# (many of the Price column values look like numbers to Excel, e.g. "592.21",
# but the source data must be preserved verbatim as text, matching the original
# file where these cells are inline strings, not numeric cells).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Auto-generated cell assignments
Set-TextValue $ws.Range("D2") "63.540.83"
$ws.Range("E2").Value = "  -0.55%  "
Set-TextValue $ws.Range("D3") "3.088.09"
$ws.Range("E3").Value = "  -1.72%  "
$ws.Range("E4").Value = "  -0.08%  "
Set-TextValue $ws.Range("D5") "592.21"
$ws.Range("E5").Value = "  +0.23%  "
Set-TextValue $ws.Range("D6") "154.99"
$ws.Range("E6").Value = "  +6.65%  "
$ws.Range("E7").Value = "  +0.03%  "
Set-TextValue $ws.Range("D8") "0.546"
$ws.Range("E8").Value = "  +3.22%  "
Set-TextValue $ws.Range("D9") "3.079.12"
$ws.Range("E9").Value = "  -1.71%  "
$ws.Range("E11").Value = "  -0.63%  "
$ws.Range("E12").Value = "  -0.30%  "
Set-TextValue $ws.Range("D13") "37.52"
$ws.Range("E13").Value = "  +0.74%  "
$ws.Range("E14").Value = "  -2.04%  "
Set-TextValue $ws.Range("D15") "3.599.70"
$ws.Range("E15").Value = "  -1.71%  "
$ws.Range("E16").Value = "  -1.45%  "
$ws.Range("E17").Value = "  -2.60%  "
Set-TextValue $ws.Range("D18") "63.517.23"
$ws.Range("E18").Value = "  -0.35%  "
Set-TextValue $ws.Range("D19") "3.084.47"
$ws.Range("E19").Value = "  -1.63%  "
Set-TextValue $ws.Range("D20") "476.60"
$ws.Range("E20").Value = "  +1.89%  "
Set-TextValue $ws.Range("D21") "14.67"
Set-TextValue $ws.Range("D22") "0.718"
$ws.Range("E22").Value = "  -1.89%  "
$ws.Range("E23").Value = "  +0.60%  "
$ws.Range("E24").Value = "  +3.99%  "
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue $ws.Range("D25") "81.27"
$ws.Range("E25").Value = "  -0.40%  "
$ws.Range("B26").Value = "InternetComputer(DFINITY)"
$ws.Range("C26").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Range("D26") "12.92"
$ws.Range("E26").Value = "  -0.64%  "
Set-TextValue $ws.Range("D27") "10.04"
$ws.Range("E27").Value = "  +2.86%  "
Set-TextValue $ws.Range("D28") "0.999"
$ws.Range("E28").Value = "  -0.14%  "
$ws.Range("E29").Value = "  -1.04%  "
$ws.Range("E30").Value = "  -0.46%  "
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("E32").Value = "  -2.25%  "
$ws.Range("E33").Value = "  +4.25%  "
Set-TextValue $ws.Range("D34") "27.27"
$ws.Range("E34").Value = "  -2.09%  "
Set-TextValue $ws.Range("D35") "0.0₃0849"
$ws.Range("E35").Value = "  +0.27%  "
Set-TextValue $ws.Range("D36") "1.05"
$ws.Range("E36").Value = "  -1.01%  "
$ws.Range("E37").Value = "  +5.63%  "
Set-TextValue $ws.Range("D38") "6.10"
$ws.Range("E38").Value = "  -0.98%  "
$ws.Range("E39").Value = "  -3.14%  "
Set-TextValue $ws.Range("D40") "9.37"
$ws.Range("E40").Value = "  -0.29%  "
$ws.Range("E41").Value = "  -1.39%  "
Set-TextValue $ws.Range("D42") "444.64"
$ws.Range("E42").Value = "  -1.93%  "
$ws.Range("E43").Value = "  -2.40%  "
Set-TextValue $ws.Range("D44") "0.0363"
$ws.Range("E44").Value = "  -2.35%  "
Set-TextValue $ws.Range("D45") "40.07"
$ws.Range("E45").Value = "  +1.19%  "
$ws.Range("E46").Value = "  +3.26%  "
Set-TextValue $ws.Range("D47") "2.801.96"
$ws.Range("E47").Value = "  -3.79%  "
Set-TextValue $ws.Range("D48") "131.54"
$ws.Range("E48").Value = "  -0.98%  "
Set-TextValue $ws.Range("D49") "25.61"
$ws.Range("E49").Value = "  +5.94%  "
$ws.Range("E50").Value = "  +0.06%  "
$ws.Range("E51").Value = "  +0.92%  "

Write-Output "Applied cryptos update"
